# Update volume tables with revised statistics/values for
# Panel B (E-mini Futures) "Ann Window Volume" (row 26),
# "Diff (Ann - Non)" (row 27) and "# Obs" (row 28) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Ann Window Volume
$ws.Range("D26").Value = 979.0687912942092
$ws.Range("E26").Value = 1215.102254259562
$ws.Range("G26").Value = 383.1935483870968
$ws.Range("H26").Value = 1558.758064516129
$ws.Range("I26").Value = 83
$ws.Range("J26").Value = 1015.054908157219
$ws.Range("K26").Value = 1239.202589920044
$ws.Range("M26").Value = 374.5245901639344
$ws.Range("N26").Value = 1716.680327868853
$ws.Range("O26").Value = 83
$ws.Range("P26").Value = 1048.295827939859
$ws.Range("Q26").Value = 1285.482856390187
$ws.Range("S26").Value = 495.8677685950413
$ws.Range("T26").Value = 1909.747933884298
$ws.Range("U26").Value = 83
$ws.Range("V26").Value = 1027.007458405049
$ws.Range("W26").Value = 1242.337598335062
$ws.Range("Y26").Value = 522.8285714285714
$ws.Range("Z26").Value = 1746.183333333333
$ws.Range("AA26").Value = 83
$ws.Range("AB26").Value = 350.7499361080687
$ws.Range("AC26").Value = 427.1553011947677
$ws.Range("AE26").Value = 289.7075757575757
$ws.Range("AF26").Value = 577.1359848484849
$ws.Range("AG26").Value = 83

# Row 27: Diff (Ann - Non)
$ws.Range("D27").Value = -10.68703847648657
$ws.Range("J27").Value = 40.69825202449143
$ws.Range("P27").Value = 89.22846758936568
$ws.Range("V27").Value = 71.2245553643144
$ws.Range("AB27").Value = 40.39626004016065

# Row 28: # Obs
$ws.Range("D28").Value = 83
$ws.Range("J28").Value = 83
$ws.Range("P28").Value = 83
$ws.Range("V28").Value = 83
$ws.Range("AB28").Value = 83
